$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2"  = "241.46"
    "D3"  = "21.32"
    "D4"  = "5.187"
    "D6"  = "3.362"
    "D8"  = "0.8022"
    "D9"  = "0.9567"
    "D10" = "0.1390"
    "D11" = "0.07326"
    "D12" = "0.03048"
    "D13" = "0.03063"
    "D14" = "0.09297"
    "D15" = "3.575"
    "D16" = "0.001641"
    "D17" = "0.04698"
    "D18" = "0.0005754"
    "D19" = "0.006435"
    "D20" = "0.004987"
    "D22" = "0.0001501"
    "D24" = "2.101"
    "D25" = "0.3241"
    "D28" = "0.0003102"
    "D40" = "0.03839"
    "D41" = "0.006890"
    "D42" = "0.1029"
    "D43" = "0.002862"
    "D44" = "0.008246"
    "D45" = "0.00005940"
    "D46" = "0.00000000750"
    "D47" = "0.0005503"
    "D48" = "0.6829"
    "D49" = "0.1148"
    "D50" = "0.00002101"
    "D51" = "0.01011"
}

foreach ($cell in $updates.Keys) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $updates[$cell]
}
